$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains weekly price observations for "Apio" at Vega Monumental
# Concepción, ordered from oldest to newest starting at row 2. A new week's
# pair of observations (Primera / Segunda quality) needs to be inserted right
# after the existing row 417, pushing every following row down by two, and
# extending the used range to R436.

# Insert two blank rows at position 418 (old rows 418.. shift down to 420..)
$ws.Rows.Item(418).Resize(2).Insert()

# New row 418: Primera quality, week of 45008
$ws.Cells.Item(418, 1).Value = 11
$ws.Cells.Item(418, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(418, 3).Value = "Bíobío"
$ws.Cells.Item(418, 4).Value = 45008
$ws.Cells.Item(418, 5).Value = 8
$ws.Cells.Item(418, 6).Value = 100112017
$ws.Cells.Item(418, 7).Value = "Apio"
$ws.Cells.Item(418, 8).Value = "Americana (o)"
$ws.Cells.Item(418, 9).Value = "Primera"
$ws.Cells.Item(418, 10).Value = 100
$ws.Cells.Item(418, 11).Value = 8000
$ws.Cells.Item(418, 12).Value = 8500
$ws.Cells.Item(418, 13).Value = 8250
$ws.Cells.Item(418, 14).Value = "`$/docena de matas"
$ws.Cells.Item(418, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(418, 16).Value = 1375
$ws.Cells.Item(418, 17).Value = 6
$ws.Cells.Item(418, 18).Value = "Hortaliza"

# New row 419: Segunda quality, same week (45008)
$ws.Cells.Item(419, 1).Value = 11
$ws.Cells.Item(419, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(419, 3).Value = "Bíobío"
$ws.Cells.Item(419, 4).Value = 45008
$ws.Cells.Item(419, 5).Value = 8
$ws.Cells.Item(419, 6).Value = 100112017
$ws.Cells.Item(419, 7).Value = "Apio"
$ws.Cells.Item(419, 8).Value = "Americana (o)"
$ws.Cells.Item(419, 9).Value = "Segunda"
$ws.Cells.Item(419, 10).Value = 50
$ws.Cells.Item(419, 11).Value = 7000
$ws.Cells.Item(419, 12).Value = 7000
$ws.Cells.Item(419, 13).Value = 7000
$ws.Cells.Item(419, 14).Value = "`$/docena de matas"
$ws.Cells.Item(419, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(419, 16).Value = 1167
$ws.Cells.Item(419, 17).Value = 6
$ws.Cells.Item(419, 18).Value = "Hortaliza"
